$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Informe Tecnic"
$ws2 = $wb.Worksheets.Item(2)   # "EXPORT_MEMORIA"

# --- Sheet1 ("Informe Tecnic") cell updates ---
$ws1.Cells.Item(351,2).Value = "23"
$ws1.Cells.Item(417,2).Value = "12613#"
$ws1.Cells.Item(417,3).Value = ""
$ws1.Cells.Item(418,2).Value = "349637#"
$ws1.Cells.Item(419,2).Value = "f2d4f20de81b017428fb8429#"
$ws1.Cells.Item(420,2).Value = "421251"
$ws1.Cells.Item(420,3).Value = "m2"
$ws1.Cells.Item(421,2).Value = "356146#"
$ws1.Cells.Item(422,2).Value = "2000920#"
$ws1.Cells.Item(423,2).Value = "360932#"
$ws1.Cells.Item(424,2).Value = "416203"
$ws1.Cells.Item(424,3).Value = "m2"
$ws1.Cells.Item(425,2).Value = "351049#"
$ws1.Cells.Item(425,3).Value = ""
$ws1.Cells.Item(426,2).Value = "2000126#"
$ws1.Cells.Item(427,2).Value = "379821"
$ws1.Cells.Item(427,3).Value = "m2"
$ws1.Cells.Item(428,2).Value = "378977#"
$ws1.Cells.Item(429,2).Value = "2001320#"
$ws1.Cells.Item(430,2).Value = "400270"
$ws1.Cells.Item(430,3).Value = "u"
$ws1.Cells.Item(431,2).Value = "01763966511d36575635ed50#"
$ws1.Cells.Item(431,3).Value = ""
$ws1.Cells.Item(432,2).Value = "2c4e32d1964278dd07b2eda4#"
$ws1.Cells.Item(433,2).Value = "68751ca122004bd195d41eb1#"
$ws1.Cells.Item(433,3).Value = ""
$ws1.Cells.Item(434,2).Value = "463aba37a91624c06d68e298#"
$ws1.Cells.Item(435,2).Value = "355123#"
$ws1.Cells.Item(436,2).Value = "2008107#"
$ws1.Cells.Item(437,2).Value = "Habit_PB#"
$ws1.Cells.Item(437,3).Value = ""
$ws1.Cells.Item(438,2).Value = "d01be03cc289732a91b81be1#"
$ws1.Cells.Item(438,3).Value = ""
$ws1.Cells.Item(439,2).Value = "384444"
$ws1.Cells.Item(439,3).Value = "m2"
$ws1.Cells.Item(440,2).Value = "3b84c223c87cf9064f14146f#"
$ws1.Cells.Item(441,2).Value = "1bdb3e9f0abc22a9cab2723d#"
$ws1.Cells.Item(442,2).Value = "54539"
$ws1.Cells.Item(442,3).Value = "m2"
$ws1.Cells.Item(443,2).Value = "8482#"
$ws1.Cells.Item(444,2).Value = "2003100#"
$ws1.Cells.Item(444,3).Value = ""
$ws1.Cells.Item(445,2).Value = "cd56c8aa009fe005b32465f1#"
$ws1.Cells.Item(446,2).Value = "2000032#"
$ws1.Cells.Item(448,2).Value = "422357"
$ws1.Cells.Item(449,2).Value = "2008163#"
$ws1.Cells.Item(450,2).Value = "Habit_P1#"
$ws1.Cells.Item(451,2).Value = "418400"
$ws1.Cells.Item(451,3).Value = "m2"
$ws1.Cells.Item(452,2).Value = "416213"
$ws1.Cells.Item(452,3).Value = "m2"
$ws1.Cells.Item(453,2).Value = "2000011#"
$ws1.Cells.Item(454,2).Value = "b74e94d46b2c057a14f47255#"
$ws1.Cells.Item(455,2).Value = "102462"
$ws1.Cells.Item(455,3).Value = "u"
$ws1.Cells.Item(457,2).Value = "395101#"
$ws1.Cells.Item(458,2).Value = "2000023#"
$ws1.Cells.Item(459,2).Value = "423688"
$ws1.Cells.Item(459,3).Value = "m2"
$ws1.Cells.Item(460,2).Value = "Cut"
$ws1.Cells.Item(460,3).Value = "m3"
$ws1.Cells.Item(461,2).Value = "2003200#"
$ws1.Cells.Item(462,2).Value = "2000051#"
$ws1.Cells.Item(463,2).Value = "141749"
$ws1.Cells.Item(463,3).Value = "u"
$ws1.Cells.Item(464,2).Value = "356145#"
$ws1.Cells.Item(464,3).Value = ""
$ws1.Cells.Item(465,2).Value = "422360"
$ws1.Cells.Item(465,3).Value = "m2"
$ws1.Cells.Item(466,2).Value = "12615"
$ws1.Cells.Item(467,2).Value = "2000120#"
$ws1.Cells.Item(468,2).Value = "385031#"
$ws1.Cells.Item(469,2).Value = "423885"
$ws1.Cells.Item(469,3).Value = "m2"
$ws1.Cells.Item(470,2).Value = "378926#"
$ws1.Cells.Item(470,3).Value = ""
$ws1.Cells.Item(471,2).Value = "88723c3731b7932f9f4568ac#"
$ws1.Cells.Item(472,2).Value = "380816"
$ws1.Cells.Item(472,3).Value = "m2"
$ws1.Cells.Item(473,2).Value = "2001330#"
$ws1.Cells.Item(474,2).Value = "421478"
$ws1.Cells.Item(474,3).Value = "m2"
$ws1.Cells.Item(475,2).Value = "Fill"
$ws1.Cells.Item(475,3).Value = "m3"
$ws1.Cells.Item(476,2).Value = "151741#"
$ws1.Cells.Item(477,2).Value = "2000160#"
$ws1.Cells.Item(478,2).Value = "Revit##0"
$ws1.Cells.Item(478,3).Value = ""
$ws1.Cells.Item(479,2).Value = "2001300#"
$ws1.Cells.Item(479,3).Value = ""
$ws1.Cells.Item(480,2).Value = "378959#"
$ws1.Cells.Item(480,3).Value = ""
$ws1.Cells.Item(481,2).Value = "49561"
$ws1.Cells.Item(481,3).Value = "u"
$ws1.Cells.Item(482,2).Value = "2000170#"
$ws1.Cells.Item(483,2).Value = "2001350#"
$ws1.Cells.Item(483,3).Value = ""
$ws1.Cells.Item(484,2).Value = "2000700#"
$ws1.Cells.Item(485,2).Value = "Áreas_P1#"
$ws1.Cells.Item(486,2).Value = "2001180#"
$ws1.Cells.Item(486,3).Value = ""
$ws1.Cells.Item(487,2).Value = "416200"
$ws1.Cells.Item(487,3).Value = "m2"
$ws1.Cells.Item(488,2).Value = "12609#"
$ws1.Cells.Item(490,2).Value = "49504"
$ws1.Cells.Item(490,3).Value = "u"
$ws1.Cells.Item(491,2).Value = "2001263#"
$ws1.Cells.Item(492,2).Value = "2000035#"
$ws1.Cells.Item(493,2).Value = "2000996#"
$ws1.Cells.Item(493,3).Value = ""
$ws1.Cells.Item(494,2).Value = "3ac1d0acb7c1b37cfcba658e#"
$ws1.Cells.Item(495,2).Value = "c5d2dd0bf170338eeffbfdb5#"
$ws1.Cells.Item(495,3).Value = ""
$ws1.Cells.Item(496,2).Value = "414195"
$ws1.Cells.Item(496,3).Value = "m2"
$ws1.Cells.Item(497,2).Value = "2003101#"
$ws1.Cells.Item(497,3).Value = ""
$ws1.Cells.Item(498,2).Value = "2000919#"
$ws1.Cells.Item(499,2).Value = "8c7cea80df2bffefc425b808#"
$ws1.Cells.Item(499,3).Value = ""
$ws1.Cells.Item(500,2).Value = "73a5c388732eabc8da8068e5#"
$ws1.Cells.Item(500,3).Value = ""
$ws1.Cells.Item(501,2).Value = "384768"
$ws1.Cells.Item(501,3).Value = "m2"
$ws1.Cells.Item(502,2).Value = "354019#"
$ws1.Cells.Item(502,3).Value = ""
$ws1.Cells.Item(503,2).Value = "368759#"
$ws1.Cells.Item(503,3).Value = ""
$ws1.Cells.Item(504,2).Value = "2000171#"
$ws1.Cells.Item(505,2).Value = "163904#"
$ws1.Cells.Item(506,2).Value = "2154006ce6da807a10c52529#"
$ws1.Cells.Item(506,3).Value = ""
$ws1.Cells.Item(507,2).Value = "416196"
$ws1.Cells.Item(507,3).Value = "m2"
$ws1.Cells.Item(508,2).Value = "395062"
$ws1.Cells.Item(508,3).Value = "u"
$ws1.Cells.Item(509,2).Value = "137323"
$ws1.Cells.Item(509,3).Value = "u"
$ws1.Cells.Item(510,2).Value = "Áreas_PB#"

# --- New row 511 (two new codes added for the flagger) ---
$ws1.Cells.Item(511,1).Value = "NO MODELAT"
$ws1.Cells.Item(511,2).Value = "ffa15592296ee41928620a0a#"
$ws1.Cells.Item(511,4).Value = "AVÍS"

# Copy the AVIS formatting (yellow fill, etc.) from the row above onto the new row
$ws1.Cells.Item(510,4).Copy()
$ws1.Cells.Item(511,4).PasteSpecial(-4122)

# --- Sheet2 ("EXPORT_MEMORIA") mirrors the same code change ---
$ws2.Cells.Item(356,2).Value = "23"
